$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.272.12"
$ws.Range("E2").Value = "'  +0.61%  "
$ws.Range("D3").Value = "'1.664.20"
$ws.Range("E3").Value = "'  +0.56%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "'  +0.76%  "
$ws.Range("D5").Value = "'218.34"
$ws.Range("E5").Value = "'  +0.22%  "
$ws.Range("E6").Value = "'  +1.35%  "
$ws.Range("D8").Value = "'0.2635"
$ws.Range("E8").Value = "'  +1.05%  "
$ws.Range("D9").Value = "'0.06360"
$ws.Range("E9").Value = "'  +0.16%  "
$ws.Range("D10").Value = "'20.52"
$ws.Range("E10").Value = "'  +0.37%  "
$ws.Range("D11").Value = "'0.07818"
$ws.Range("E11").Value = "'  +0.44%  "
$ws.Range("D12").Value = "'4.565"
$ws.Range("E12").Value = "'  +1.45%  "
$ws.Range("D13").Value = "'1.671.08"
$ws.Range("E13").Value = "'  +3.20%  "
$ws.Range("D14").Value = "'1.891.60"
$ws.Range("E14").Value = "'  +0.44%  "
$ws.Range("D15").Value = "'0.5533"
$ws.Range("E15").Value = "'  +0.92%  "
$ws.Range("D16").Value = "'0.0₅8194"
$ws.Range("E16").Value = "'  -0.50%  "
$ws.Range("D17").Value = "'65.66"
$ws.Range("E17").Value = "'  +0.43%  "
$ws.Range("E18").Value = "'  +0.78%  "
$ws.Range("D19").Value = "'4.682"
$ws.Range("E19").Value = "'  +2.24%  "
$ws.Range("D20").Value = "'193.96"
$ws.Range("E20").Value = "'  +1.19%  "
$ws.Range("D21").Value = "'10.19"
$ws.Range("E21").Value = "'  +1.18%  "
$ws.Range("D22").Value = "'6.035"
$ws.Range("E22").Value = "'  -0.07%  "
$ws.Range("E23").Value = "'  +0.74%  "
$ws.Range("D24").Value = "'146.02"
$ws.Range("E24").Value = "'  +2.85%  "
$ws.Range("D25").Value = "'0.1226"
$ws.Range("E25").Value = "'  -1.80%  "
$ws.Range("D26").Value = "'7.191"
$ws.Range("E26").Value = "'  -1.03%  "
$ws.Range("D27").Value = "'16.14"
$ws.Range("E27").Value = "'  +0.07%  "
$ws.Range("D28").Value = "'1.488"
$ws.Range("E28").Value = "'  +3.87%  "
$ws.Range("D29").Value = "'0.05869"
$ws.Range("E29").Value = "'  -0.66%  "
$ws.Range("D30").Value = "'1.279"
$ws.Range("E30").Value = "'  -0.23%  "
$ws.Range("D31").Value = "'3.586"
$ws.Range("E31").Value = "'  +1.73%  "
$ws.Range("D32").Value = "'3.279"
$ws.Range("E32").Value = "'  +0.73%  "
$ws.Range("E33").Value = "'  +1.41%  "
$ws.Range("D34").Value = "'0.9609"
$ws.Range("E34").Value = "'  +0.82%  "
$ws.Range("D35").Value = "'2.822"
$ws.Range("E35").Value = "'  +1.44%  "
$ws.Range("E36").Value = "'  +0.54%  "
$ws.Range("D37").Value = "'0.5800"
$ws.Range("E37").Value = "'  +1.71%  "
$ws.Range("D38").Value = "'0.01604"
$ws.Range("E38").Value = "'  -0.93%  "
$ws.Range("D39").Value = "'0.8640"
$ws.Range("E39").Value = "'  +1.86%  "
$ws.Range("D40").Value = "'5.839"
$ws.Range("E40").Value = "'  +0.94%  "
$ws.Range("E41").Value = "'  +0.73%  "
$ws.Range("D42").Value = "'1.047.33"
$ws.Range("E42").Value = "'  +1.79%  "
$ws.Range("D43").Value = "'104.05"
$ws.Range("E43").Value = "'  +0.91%  "
$ws.Range("D44").Value = "'1.802.33"
$ws.Range("E44").Value = "'  +0.22%  "
$ws.Range("D45").Value = "'57.61"
$ws.Range("E45").Value = "'  +0.58%  "
$ws.Range("E46").Value = "'  +0.97%  "
$ws.Range("E47").Value = "'  -5.48%  "
$ws.Range("D48").Value = "'0.4380"
$ws.Range("E48").Value = "'  +1.85%  "
$ws.Range("D49").Value = "'8.052"
$ws.Range("E49").Value = "'  +2.37%  "
$ws.Range("D50").Value = "'0.05161"
$ws.Range("D51").Value = "'1.428"
$ws.Range("E51").Value = "'  -3.36%  "
